# Update PatientData.xlsx with latest test data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update patient name
$ws.Range("A2").Value = "Olivia"
$ws.Range("B2").Value = "Kerr"

# Update phone number
$ws.Range("D2").Value = 3558693829

# Update date of birth (stored as text)
$ws.Range("E2").Value = "02/09/2001"

# Update email + hyperlink (C2). Adding a hyperlink via COM auto-applies
# Excel's built-in "Hyperlink" style, but this cell already carries its own
# explicit blue-font formatting, so capture it first and restore it after
# the hyperlink is rebuilt to point at the new address/display text.
$c = $ws.Range("C2")

$fontColor = $c.Font.Color
$fontName = $c.Font.Name
$fontSize = $c.Font.Size
$fontUnderline = $c.Font.Underline
$halign = $c.HorizontalAlignment
$valign = $c.VerticalAlignment
$wrap = $c.WrapText
$numfmt = $c.NumberFormat

$c.Value = "olivia@k.com"
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($c, "mailto:olivia@k.com", "", "", "olivia@k.com")

$c.Font.Color = $fontColor
$c.Font.Name = $fontName
$c.Font.Size = $fontSize
$c.Font.Underline = $fontUnderline
$c.HorizontalAlignment = $halign
$c.VerticalAlignment = $valign
$c.WrapText = $wrap
$c.NumberFormat = $numfmt
